# Insert a new daily price record (row 277) for "Navel Late" into the
# weekly/daily Naranja - Vega Monumental Concepción sheet.
# All existing rows from 277 downward shift down by one (277->278, ..., 332->333).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 277:332 down to 278:333 by inserting a fresh row at 277.
$ws.Rows("277:277").Insert()

# Populate the newly inserted row 277 with the new record.
$ws.Range("A277").Value = 11
$ws.Range("B277").Value = "Vega Monumental Concepción"
$ws.Range("C277").Value = "Bíobío"
$ws.Range("D277").Value = 44798
$ws.Range("E277").Value = 8
$ws.Range("F277").Value = "Fruta"
$ws.Range("G277").Value = 100102
$ws.Range("H277").Value = "Cítricos"
$ws.Range("I277").Value = 100102005
$ws.Range("J277").Value = "Naranja"
$ws.Range("K277").Value = "Navel Late"
$ws.Range("L277").Value = "Primera"
$ws.Range("M277").Value = 350
$ws.Range("N277").Value = 5500
$ws.Range("O277").Value = 6000
$ws.Range("P277").Value = 5714
$ws.Range("Q277").Value = "$/caja 15 kilos"
$ws.Range("R277").Value = "Región de O'Higgins"
$ws.Range("S277").Value = 381
$ws.Range("T277").Value = 15
